$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dataRange = $ws.Range("A2:B57")
$keyRange = $ws.Range("B2:B57")

$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($keyRange, 0, 1, $null, 0) | Out-Null

$ws.Sort.SetRange($dataRange)
$ws.Sort.Header = 2
$ws.Sort.MatchCase = $false
$ws.Sort.Orientation = 1
$ws.Sort.SortMethod = 1
$ws.Sort.Apply()
